$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Ntf3"
$ws.Cells.Item(2, 3).Value = "Ngfr"
$ws.Cells.Item(2, 4).Value = "FAPs"
$ws.Cells.Item(2, 5).Value = 2
$ws.Cells.Item(2, 6).Value = 0.6666666666666666
$ws.Cells.Item(2, 7).Value = 4.477828666666666
$ws.Cells.Item(2, 8).Value = 13.433486
$ws.Cells.Item(2, 9).Value = 0.4652827882180238
$ws.Cells.Item(2, 10).Value = 0.4652827882180238
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 0.6946430000000001
$ws.Cells.Item(2, 14).Value = 2.083929
$ws.Cells.Item(2, 15).Value = 0.1269399741689062
$ws.Cells.Item(2, 16).Value = 0.1269399741689062
$ws.Cells.Item(2, 17).Value = 3.110492338499334
$ws.Cells.Item(2, 18).Value = 27.994431046494
$ws.Cells.Item(2, 19).Value = 0.0590629851176326
$ws.Cells.Item(2, 20).Value = 0.0590629851176326

# Row 3
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Ntf3"
$ws.Cells.Item(3, 3).Value = "Ngfr"
$ws.Cells.Item(3, 4).Value = "M2"
$ws.Cells.Item(3, 5).Value = 2
$ws.Cells.Item(3, 6).Value = 0.6666666666666666
$ws.Cells.Item(3, 7).Value = 4.477828666666666
$ws.Cells.Item(3, 8).Value = 13.433486
$ws.Cells.Item(3, 9).Value = 0.4652827882180238
$ws.Cells.Item(3, 10).Value = 0.4652827882180238
$ws.Cells.Item(3, 11).Value = 1
$ws.Cells.Item(3, 12).Value = 0.3333333333333333
$ws.Cells.Item(3, 13).Value = 0.05042666666666667
$ws.Cells.Item(3, 14).Value = 0.15128
$ws.Cells.Item(3, 15).Value = 0.009215035297398391
$ws.Cells.Item(3, 16).Value = 0.009215035297398391
$ws.Cells.Item(3, 17).Value = 0.2258019735644445
$ws.Cells.Item(3, 18).Value = 2.03221776208
$ws.Cells.Item(3, 19).Value = 0.004287597316701029
$ws.Cells.Item(3, 20).Value = 0.004287597316701029

# Row 4
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Ntf3"
$ws.Cells.Item(4, 3).Value = "Ngfr"
$ws.Cells.Item(4, 4).Value = "sCs"
$ws.Cells.Item(4, 5).Value = 2
$ws.Cells.Item(4, 6).Value = 0.6666666666666666
$ws.Cells.Item(4, 7).Value = 4.477828666666666
$ws.Cells.Item(4, 8).Value = 13.433486
$ws.Cells.Item(4, 9).Value = 0.4652827882180238
$ws.Cells.Item(4, 10).Value = 0.4652827882180238
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 4.727146666666667
$ws.Cells.Item(4, 14).Value = 14.18144
$ws.Cells.Item(4, 15).Value = 0.8638449905336953
$ws.Cells.Item(4, 16).Value = 0.8638449905336955
$ws.Cells.Item(4, 17).Value = 21.16735285553778
$ws.Cells.Item(4, 18).Value = 190.50617569984
$ws.Cells.Item(4, 19).Value = 0.4019322057836901
$ws.Cells.Item(4, 20).Value = 0.4019322057836902

# Row 5
$ws.Cells.Item(5, 1).Value = "FAPs"
$ws.Cells.Item(5, 2).Value = "Ntf3"
$ws.Cells.Item(5, 3).Value = "Ngfr"
$ws.Cells.Item(5, 4).Value = "FAPs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 2.528563666666667
$ws.Cells.Item(5, 8).Value = 7.585691000000001
$ws.Cells.Item(5, 9).Value = 0.2627383137214249
$ws.Cells.Item(5, 10).Value = 0.2627383137214249
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 0.6946430000000001
$ws.Cells.Item(5, 14).Value = 2.083929
$ws.Cells.Item(5, 15).Value = 0.1269399741689062
$ws.Cells.Item(5, 16).Value = 0.1269399741689062
$ws.Cells.Item(5, 17).Value = 1.756449051104334
$ws.Cells.Item(5, 18).Value = 15.808041459939
$ws.Cells.Item(5, 19).Value = 0.03335199475697966
$ws.Cells.Item(5, 20).Value = 0.03335199475697966

# Row 6
$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Ntf3"
$ws.Cells.Item(6, 3).Value = "Ngfr"
$ws.Cells.Item(6, 4).Value = "M2"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 2.528563666666667
$ws.Cells.Item(6, 8).Value = 7.585691000000001
$ws.Cells.Item(6, 9).Value = 0.2627383137214249
$ws.Cells.Item(6, 10).Value = 0.2627383137214249
$ws.Cells.Item(6, 11).Value = 1
$ws.Cells.Item(6, 12).Value = 0.3333333333333333
$ws.Cells.Item(6, 13).Value = 0.05042666666666667
$ws.Cells.Item(6, 14).Value = 0.15128
$ws.Cells.Item(6, 15).Value = 0.009215035297398391
$ws.Cells.Item(6, 16).Value = 0.009215035297398391
$ws.Cells.Item(6, 17).Value = 0.1275070371644444
$ws.Cells.Item(6, 18).Value = 1.14756333448
$ws.Cells.Item(6, 19).Value = 0.002421142834921862
$ws.Cells.Item(6, 20).Value = 0.002421142834921862

# Row 7
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Ntf3"
$ws.Cells.Item(7, 3).Value = "Ngfr"
$ws.Cells.Item(7, 4).Value = "sCs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 2.528563666666667
$ws.Cells.Item(7, 8).Value = 7.585691000000001
$ws.Cells.Item(7, 9).Value = 0.2627383137214249
$ws.Cells.Item(7, 10).Value = 0.2627383137214249
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 4.727146666666667
$ws.Cells.Item(7, 14).Value = 14.18144
$ws.Cells.Item(7, 15).Value = 0.8638449905336953
$ws.Cells.Item(7, 16).Value = 0.8638449905336955
$ws.Cells.Item(7, 17).Value = 11.95289130833778
$ws.Cells.Item(7, 18).Value = 107.57602177504
$ws.Cells.Item(7, 19).Value = 0.2269651761295234
$ws.Cells.Item(7, 20).Value = 0.2269651761295234

# Row 8
$ws.Cells.Item(8, 1).Value = "sCs"
$ws.Cells.Item(8, 2).Value = "Ntf3"
$ws.Cells.Item(8, 3).Value = "Ngfr"
$ws.Cells.Item(8, 4).Value = "FAPs"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 2.617494
$ws.Cells.Item(8, 8).Value = 7.852482
$ws.Cells.Item(8, 9).Value = 0.2719788980605514
$ws.Cells.Item(8, 10).Value = 0.2719788980605514
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 0.6946430000000001
$ws.Cells.Item(8, 14).Value = 2.083929
$ws.Cells.Item(8, 15).Value = 0.1269399741689062
$ws.Cells.Item(8, 16).Value = 0.1269399741689062
$ws.Cells.Item(8, 17).Value = 1.818223884642
$ws.Cells.Item(8, 18).Value = 16.364014961778
$ws.Cells.Item(8, 19).Value = 0.03452499429429397
$ws.Cells.Item(8, 20).Value = 0.03452499429429397

# Row 9
$ws.Cells.Item(9, 1).Value = "sCs"
$ws.Cells.Item(9, 2).Value = "Ntf3"
$ws.Cells.Item(9, 3).Value = "Ngfr"
$ws.Cells.Item(9, 4).Value = "M2"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 2.617494
$ws.Cells.Item(9, 8).Value = 7.852482
$ws.Cells.Item(9, 9).Value = 0.2719788980605514
$ws.Cells.Item(9, 10).Value = 0.2719788980605514
$ws.Cells.Item(9, 11).Value = 1
$ws.Cells.Item(9, 12).Value = 0.3333333333333333
$ws.Cells.Item(9, 13).Value = 0.05042666666666667
$ws.Cells.Item(9, 14).Value = 0.15128
$ws.Cells.Item(9, 15).Value = 0.009215035297398391
$ws.Cells.Item(9, 16).Value = 0.009215035297398391
$ws.Cells.Item(9, 17).Value = 0.13199149744
$ws.Cells.Item(9, 18).Value = 1.18792347696
$ws.Cells.Item(9, 19).Value = 0.0025062951457755
$ws.Cells.Item(9, 20).Value = 0.0025062951457755

# Row 10
$ws.Cells.Item(10, 1).Value = "sCs"
$ws.Cells.Item(10, 2).Value = "Ntf3"
$ws.Cells.Item(10, 3).Value = "Ngfr"
$ws.Cells.Item(10, 4).Value = "sCs"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 2.617494
$ws.Cells.Item(10, 8).Value = 7.852482
$ws.Cells.Item(10, 9).Value = 0.2719788980605514
$ws.Cells.Item(10, 10).Value = 0.2719788980605514
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 4.727146666666667
$ws.Cells.Item(10, 14).Value = 14.18144
$ws.Cells.Item(10, 15).Value = 0.8638449905336953
$ws.Cells.Item(10, 16).Value = 0.8638449905336955
$ws.Cells.Item(10, 17).Value = 12.37327803712
$ws.Cells.Item(10, 18).Value = 111.35950233408
$ws.Cells.Item(10, 19).Value = 0.2349476086204819
$ws.Cells.Item(10, 20).Value = 0.2349476086204819
